# =====================================================================
# issue #5: property land done
#
# Sheet "land" (tab 1, Chinese name stored in the sheet itself): the
# Chinese column headers (B1:H1) are replaced with the canonical English
# field names used by the scraping pipeline, and 7 metadata columns are
# appended after them (I:O): property_category, category, date,
# legislator_name, legislator_id, source_file, index - populated for
# every one of the 4 data rows. A handful of text values are also
# cleaned up across this and two other sheets (stray spaces / hyphens
# removed from parcel numbers and dates, a stray trailing dot removed).
# =====================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1 (tab 1): land
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Give the 7 new header cells (I1:O1) the same bold/bordered/centered
# format already used by the rest of row 1 (format-only copy/paste).
$ws1.Range("H1").Copy() | Out-Null
$ws1.Range("I1:O1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Row 1 header labels
$ws1.Range("B1").Value = "name"
$ws1.Range("C1").Value = "area"
$ws1.Range("D1").Value = "share_portion"
$ws1.Range("E1").Value = "owner"
$ws1.Range("F1").Value = "register_date"
$ws1.Range("G1").Value = "register_reason"
$ws1.Range("H1").Value = "acquire_value"
$ws1.Range("I1").Value = "property_category"
$ws1.Range("J1").Value = "category"
$ws1.Range("K1").Value = "date"
$ws1.Range("L1").Value = "legislator_name"
$ws1.Range("M1").Value = "legislator_id"
$ws1.Range("N1").Value = "source_file"
$ws1.Range("O1").Value = "index"

# Row 2
$ws1.Range("B2").Value = "新北市三峽區十三添段十三添小段00050027地號"
$ws1.Range("F2").Value = "93年04月08日"
$ws1.Range("I2").Value = "land"
$ws1.Range("J2").Value = "normal"
$ws1.Range("K2").Value = "2012-04-30"
$ws1.Range("L2").Value = "羅明才"
$ws1.Range("M2").Value = 879
$ws1.Range("N2").Value = "tmpa5201"
$ws1.Range("O2").Value = 15

# Row 3
$ws1.Range("B3").Value = "臺北市中山區吉林段四小段08840000地號"
$ws1.Range("D3").Value = "10000分之44"
$ws1.Range("F3").Value = "75年12月04日"
$ws1.Range("G3").Value = "買賣"
$ws1.Range("I3").Value = "land"
$ws1.Range("J3").Value = "normal"
$ws1.Range("K3").Value = "2012-04-30"
$ws1.Range("L3").Value = "羅明才"
$ws1.Range("M3").Value = 879
$ws1.Range("N3").Value = "tmpa5201"
$ws1.Range("O3").Value = 16

# Row 4
$ws1.Range("B4").Value = "臺北市中山區吉林段四小段08890000地號"
$ws1.Range("D4").Value = "10000分之44"
$ws1.Range("F4").Value = "75年12月04日"
$ws1.Range("I4").Value = "land"
$ws1.Range("J4").Value = "normal"
$ws1.Range("K4").Value = "2012-04-30"
$ws1.Range("L4").Value = "羅明才"
$ws1.Range("M4").Value = 879
$ws1.Range("N4").Value = "tmpa5201"
$ws1.Range("O4").Value = 17

# Row 5
$ws1.Range("B5").Value = "臺北市中山區吉林段四小段08890001地號"
$ws1.Range("D5").Value = "10000分之44"
$ws1.Range("F5").Value = "75年12月04日"
$ws1.Range("I5").Value = "land"
$ws1.Range("J5").Value = "normal"
$ws1.Range("K5").Value = "2012-04-30"
$ws1.Range("L5").Value = "羅明才"
$ws1.Range("M5").Value = 879
$ws1.Range("N5").Value = "tmpa5201"
$ws1.Range("O5").Value = 18

# ---------------------------------------------------------------------
# Sheet 2 (tab 2): building - text clean-up only, no structural change
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "臺北市中山區吉林段四小段00638000建號"
$ws2.Range("F2").Value = "75年12月04日"
$ws2.Range("B3").Value = "新北市三峽區十三添段十三添小段00091000建號"
$ws2.Range("F3").Value = "93年04月08日"

# ---------------------------------------------------------------------
# Sheet 3 (tab 3): deposit - text clean-up only, no structural change
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B4").Value = "中國信託商業銀行板橋分行"
$ws3.Range("B9").Value = "中國信託商業銀行板橋分行"
